# Doc/Wykonane zadania.xlsx -- "+/- instruktora do pośredniej kolekcji"
#
# Adds two new task rows (A6:A7) with completion dates in column B,
# a trailing date-only entry in B8, and backfills a date in column B
# for the existing rows (B2:B5). Also adds the TRUE/TRUNC helper
# formula pair (O15:O16) used elsewhere for date-processing logic,
# and leaves the sheet dimension/selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- new task descriptions (new shared strings) -------------------------
$ws.Range("A6").Value = "funkcja vba konwertująca polskie znaki na łacińskie"
$ws.Range("A7").Value = "implementacja linii poleceń w LanguageCourseView"

# --- completion dates -----------------------------------------------
# serial 41901 = 2014-09-19, 41904 = 2014-09-22, 41905 = 2014-09-23
$ws.Range("B2").Value = 41901
$ws.Range("B3").Value = 41901
$ws.Range("B4").Value = 41901
$ws.Range("B5").Value = 41901
$ws.Range("B6").Value = 41904
$ws.Range("B7").Value = 41904
$ws.Range("B8").Value = 41905

# apply the built-in short-date format (numFmtId 14) to B2, then fan it
# out to the rest of the column via copy/paste-format so every cell
# shares the same style record instead of minting a new one each time
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy()
$ws.Range("B3:B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# size column B to fit the new date values
$ws.Columns.Item(2).AutoFit()

# --- helper formulas used elsewhere for date-processing view logic -----
$ws.Range("O15").Formula = "=TRUE"
$ws.Range("O16").Formula = "=TRUNC(O15)"

# --- restore the cursor to where the author left it --------------------
$ws.Range("B26").Select()
